# Daily attendance processing - 2025-10-22 11:42:44
# Swap the order of "System" and the recorder's email in the
# "Recorded By" column (G) wherever the value is exactly
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Text -eq $target) {
        $cell.Value = $replacement
    }
}
